$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12/13 swap: Polkadot <-> TRON (name, link, price, volume) ---
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.114"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.15%  "

# --- Price / Volume(1h) updates for all other rows ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.372.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.525.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.72%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.912.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.96%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.523.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.476.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("E24").Value = "  +1.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "

# Row 28
$ws.Range("E28").Value = "  +1.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "

# Row 31
$ws.Range("E31").Value = "  -4.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "

# Row 38
$ws.Range("E38").Value = "  -3.59%  "

# Row 39
$ws.Range("E39").Value = "  -1.41%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.09%  "

# Row 42
$ws.Range("E42").Value = "  -2.66%  "

# Row 43
$ws.Range("E43").Value = "  +0.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0298"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.026.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.64%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.767.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
